$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fix the header formula row (row 21) - column C changed to a Contains()-based condition
$ws.Range("C21").Value = '(((CaseFilePipelineContext)$model.getPipelineContext()).getEnqueueName() != null) && ((CaseFilePipelineContext)$model.getPipelineContext()).getEnqueueName().contains("$param")'

# Remove old rule-table rows (23-43) entirely, shifting styles out with them,
# so the freshly typed replacement rows start from a clean/default style.
$ws.Range("A23:G43").EntireRow.Delete()

# Re-populate the rule table with the restructured FOIA rules
$ws.Range("B23").Value = 'request type'
$ws.Range("C23").Value = 'Fulfill,Approve,General Counsel,Billing,Release'
$ws.Range("E23").Value = 'requestType == null || requestType?.trim().isEmpty()'
$ws.Range("F23").Value = 'Request type is required'

$ws.Range("B24").Value = 'expedite flag'
$ws.Range("C24").Value = 'Fulfill,Approve,General Counsel,Billing,Release'
$ws.Range("E24").Value = 'expediteFlag == null'
$ws.Range("F24").Value = 'Expedite flag is required'

$ws.Range("B25").Value = 'fee waiver flag'
$ws.Range("C25").Value = 'Fulfill,Approve,General Counsel,Billing,Release'
$ws.Range("E25").Value = 'feeWaiverFlag == null'
$ws.Range("F25").Value = 'Fee waiver flag is required'

$ws.Range("B26").Value = 'litigation flag'
$ws.Range("C26").Value = 'Fulfill,Approve,General Counsel,Billing,Release'
$ws.Range("E26").Value = 'litigationFlag == null'
$ws.Range("F26").Value = 'Litigation flag is required'

$ws.Range("B27").Value = 'Must be litigation'
$ws.Range("C27").Value = 'General Counsel'
$ws.Range("E27").Value = 'litigationFlag == null || !litigationFlag'
$ws.Range("F27").Value = 'Litigation flag must be checked'

$ws.Range("B28").Value = 'Must not have waived fees'
$ws.Range("C28").Value = 'Billing'
$ws.Range("E28").Value = 'feeWaiverFlag != null && feeWaiverFlag'
$ws.Range("F28").Value = 'The fee waiver flag must not be checked'

$ws.Range("B29").Value = 'Must have disposition type'
$ws.Range("C29").Value = 'Approve'
$ws.Range("E29").Value = 'disposition == null && deniedFlag'
$ws.Range("F29").Value = 'Disposition Type is required'

$ws.Range("B30").Value = 'Must have disposition subtype (when disposition type is full-denial)'
$ws.Range("C30").Value = 'Approve'
$ws.Range("E30").Value = 'disposition != null && disposition.equals("full-denial") && dispositionSubtype == null && deniedFlag'
$ws.Range("F30").Value = 'Disposition SubType is required'

$ws.Range("B31").Value = 'Must have other field (when disposition is full denial and subtype is other)'
$ws.Range("C31").Value = 'Approve'
$ws.Range("E31").Value = 'dispositionSubtype == ''other'' && otherReason == null && deniedFlag'
$ws.Range("F31").Value = 'Disposition "Other" reason is required'

$ws.Range("B32").Value = 'Must have disposition type'
$ws.Range("C32").Value = 'Release'
$ws.Range("D32").Value = 'Approve'
$ws.Range("E32").Value = 'disposition == null && deniedFlag'
$ws.Range("F32").Value = 'Disposition Type is required'

$ws.Range("B33").Value = 'Must have disposition subtype (when disposition type is full-denial)'
$ws.Range("C33").Value = 'Release'
$ws.Range("D33").Value = 'Approve'
$ws.Range("E33").Value = 'disposition != null && disposition.equals("full-denial") && dispositionSubtype == null && deniedFlag'
$ws.Range("F33").Value = 'Disposition SubType is required'

$ws.Range("B34").Value = 'Must have other field (when disposition is full denial and subtype is other)'
$ws.Range("C34").Value = 'Release'
$ws.Range("D34").Value = 'Approve'
$ws.Range("E34").Value = 'dispositionSubtype == ''other'' && otherReason == null && deniedFlag'
$ws.Range("F34").Value = 'Disposition "Other" reason is required'

$ws.Range("B35").Value = 'Must have set Executive Group'
$ws.Range("C35").Value = 'Fulfill'
$ws.Range("E35").Value = 'notificationGroup == null'
$ws.Range("F35").Value = 'Executive Group is required'

$ws.Range("B36").Value = 'Must have received date'
$ws.Range("C36").Value = 'Fulfill'
$ws.Range("D36").Value = 'Intake'
$ws.Range("E36").Value = 'receivedDate==null'
$ws.Range("F36").Value = 'Received date is required'

# Update the view so the selection/scroll position matches the trimmed table
$ws.Range("A36").Select()
$excel.ActiveWindow.ScrollRow = 20
